# ADMS Offer Letter - PDF/import finish
# - Row 2 now holds the real imported record (Employee ID / Client Company / Offer format)
# - Rows 3-7 (the old sample/test rows) are wiped back to blank, templated rows
# - Column B is widened and the "Offer Letter For" cell gets a highlighted font

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 2: the single real data row -------------------------------------------------
$ws.Range("A2").Value = "FFF1616"
$ws.Range("B2").Value = "Digital Age Retail Private Limited "
$ws.Range("C2").Value = "Format 1"

$ws.Rows.Item(2).RowHeight = 18.75

# Highlight the "Offer Letter For" cell
$ws.Range("C2").Font.Name = "Consolas"
$ws.Range("C2").Font.Size = 14
$ws.Range("C2").Font.Color = 7658470

# --- Rows 3-7: drop the old sample rows, leave a formatted blank template row behind ---
$ws.Range("A3:BQ7").ClearContents()

$ws.Range("B3").Font.Name = "Arial"
$ws.Range("B3").Font.Size = 10
$ws.Range("B3").Font.Color = 0

$ws.Range("B4").Font.Name = "Arial"
$ws.Range("B4").Font.Size = 10
$ws.Range("B4").Font.Color = 0

$ws.Range("B5").Font.Name = "Arial"
$ws.Range("B5").Font.Size = 10
$ws.Range("B5").Font.Color = 0

$ws.Range("B6").Font.Name = "Arial"
$ws.Range("B6").Font.Size = 10
$ws.Range("B6").Font.Color = 0

$ws.Range("B7").Font.Name = "Arial"
$ws.Range("B7").Font.Size = 10
$ws.Range("B7").Font.Color = 0

# --- Layout tweaks ----------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 38.5

# --- Selection ----------------------------------------------------------------------------
$ws.Range("C2").Select()
